# Applies the stimuli-sheet updates described in the commit diff:
#  - adds a "carrier" value (D column) to the practice rows (2-5), matching
#    the pair_kind already recorded in column K for that row
#  - adds a "unique_video"/"unique_audio" marker (J column) to rows 6-9
#  - populates the previously-empty rows 14-21 with the matching C
#    (unique_video/unique_audio) and D (can/where/do/look) values, mirroring
#    rows 6-13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows 2-5: new D column = the carrier word used for that practice pair
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Rows 6-9: add J column = unique_video / unique_audio
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: populate C (unique_video/unique_audio) and D (can/where/do/look),
# mirroring the existing carrier pattern from rows 6-13
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"
$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
